# "scenario to be ran" - update the Coupling Parameters sheet with the new
# scenario values, drop the two formulas that were replaced by their plain
# results, and move the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# End Year: 2030 -> 2050
$ws.Range("B4").Value = 2050

# Look Ahead: was "=B9" (evaluates to 4) -> plain static value 4
$ws.Range("B8").Value = 4

# short_term_investment_minimal_irr: 0.3 -> 0.2
$ws.Range("B10").Value = 0.2

# InvestmentIteration: 10 -> 40
$ws.Range("B11").Value = 40

# max_permit_build_time: 10 -> 4
$ws.Range("B12").Value = 4

# start_year_dismantling: 150 -> 10000
$ws.Range("B13").Value = 10000

# typeofProfitforPastHorizon: "none" -> "totalProfits"
$ws.Range("B14").Value = "totalProfits"

# fix_prices_to_2020: was "=IF(B13>=10000,FALSE,TRUE)" -> plain static TRUE
$ws.Range("B16").Value = $true

# Move the active selection from C8 to C6
$ws.Range("C6").Select()
